$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NDBKY")

$ws.Range("D8").Value = 2585100
$ws.Range("E8").Value = 2600100
$ws.Range("F8").Value = 2560900
$ws.Range("G8").Value = 2754100
$ws.Range("H8").Value = 2276400
$ws.Range("I8").Value = 2177900
$ws.Range("J8").Value = 1954300
$ws.Range("K8").Value = 1936000
$ws.Range("D15").Value = -79100
$ws.Range("E15").Value = -81200
$ws.Range("F15").Value = -75600
$ws.Range("G15").Value = -75700
$ws.Range("H15").Value = -69500
$ws.Range("I15").Value = -64800
$ws.Range("J15").Value = -59300
$ws.Range("K15").Value = -61300
$ws.Range("D17").Value = 1749500
$ws.Range("E17").Value = 1752500
$ws.Range("F17").Value = 1741600
$ws.Range("G17").Value = 1996400
$ws.Range("H17").Value = 1535000
$ws.Range("I17").Value = 1511200
$ws.Range("J17").Value = 1312200
$ws.Range("K17").Value = 1261400
$ws.Range("D18").Value = 835600
$ws.Range("E18").Value = 847600
$ws.Range("F18").Value = 819300
$ws.Range("G18").Value = 757700
$ws.Range("H18").Value = 741400
$ws.Range("I18").Value = 666800
$ws.Range("J18").Value = 642100
$ws.Range("K18").Value = 674600
$ws.Range("D20").Value = -194900
$ws.Range("E20").Value = -247800
$ws.Range("F20").Value = -287700
$ws.Range("G20").Value = -276600
$ws.Range("H20").Value = -220800
$ws.Range("I20").Value = -164500
$ws.Range("J20").Value = -138100
$ws.Range("K20").Value = -161300
$ws.Range("D21").Value = 719800
$ws.Range("E21").Value = 605400
$ws.Range("F21").Value = 607300
$ws.Range("I21").Value = 507700
$ws.Range("J21").Value = 563300
$ws.Range("K21").Value = 574600
$ws.Range("D23").Value = 640700
$ws.Range("E23").Value = 599800
$ws.Range("F23").Value = 531700
$ws.Range("G23").Value = 481100
$ws.Range("H23").Value = 520600
$ws.Range("I23").Value = 502300
$ws.Range("J23").Value = 504000
$ws.Range("K23").Value = 513300
$ws.Range("D24").Value = 160800
$ws.Range("E24").Value = 139200
$ws.Range("F24").Value = 149300
$ws.Range("G24").Value = 137800
$ws.Range("H24").Value = 133200
$ws.Range("I24").Value = 116400
$ws.Range("J24").Value = 124700
$ws.Range("K24").Value = 130400
$ws.Range("D26").Value = 479900
$ws.Range("E26").Value = 460600
$ws.Range("F26").Value = 382400
$ws.Range("G26").Value = 343200
$ws.Range("H26").Value = 387300
$ws.Range("I26").Value = 385800
$ws.Range("J26").Value = 379200
$ws.Range("K26").Value = 382900
$ws.Range("D27").Value = 456300
$ws.Range("E27").Value = 437100
$ws.Range("F27").Value = 359400
$ws.Range("G27").Value = 321500
$ws.Range("H27").Value = 373000
$ws.Range("I27").Value = 369600
$ws.Range("J27").Value = 365200
$ws.Range("K27").Value = 368100
$ws.Range("D32").Value = 194900
$ws.Range("E32").Value = 247800
$ws.Range("F32").Value = 287700
$ws.Range("G32").Value = 276600
$ws.Range("H32").Value = 220800
$ws.Range("I32").Value = 164500
$ws.Range("J32").Value = 138100
$ws.Range("K32").Value = 161300
$ws.Range("D33").Value = 456300
$ws.Range("E33").Value = 437100
$ws.Range("F33").Value = 359400
$ws.Range("G33").Value = 321500
$ws.Range("H33").Value = 373000
$ws.Range("I33").Value = 369600
$ws.Range("J33").Value = 365200
$ws.Range("K33").Value = 368100
$ws.Range("D35").Value = 456300
$ws.Range("E35").Value = 437100
$ws.Range("F35").Value = 359400
$ws.Range("G35").Value = 321500
$ws.Range("H35").Value = 373000
$ws.Range("I35").Value = 369600
$ws.Range("J35").Value = 365200
$ws.Range("K35").Value = 368100
$ws.Range("D41").Value = 3198500
$ws.Range("E41").Value = 3880500
$ws.Range("F41").Value = 4303600
$ws.Range("G41").Value = 5117300
$ws.Range("H41").Value = 4749200
$ws.Range("I41").Value = 4715600
$ws.Range("J41").Value = 5383100
$ws.Range("K41").Value = 3554500
$ws.Range("D42").Value = 6458300
$ws.Range("E42").Value = 6358800
$ws.Range("F42").Value = 6219400
$ws.Range("G42").Value = 6717300
$ws.Range("H42").Value = 6470400
$ws.Range("I42").Value = 6051400
$ws.Range("J42").Value = 5828200
$ws.Range("K42").Value = 6145400
$ws.Range("D47").Value = 212700
$ws.Range("E47").Value = 460700
$ws.Range("F47").Value = 408800
$ws.Range("G47").Value = 450100
$ws.Range("H47").Value = 546100
$ws.Range("I47").Value = 656500
$ws.Range("J47").Value = 489800
$ws.Range("K47").Value = 543200
$ws.Range("D48").Value = 591500
$ws.Range("E48").Value = 610100
$ws.Range("F48").Value = 593900
$ws.Range("G48").Value = 616200
$ws.Range("H48").Value = 612700
$ws.Range("I48").Value = 604200
$ws.Range("J48").Value = 536400
$ws.Range("K48").Value = 559700
$ws.Range("D49").Value = 835100
$ws.Range("E49").Value = 780300
$ws.Range("F49").Value = 723900
$ws.Range("G49").Value = 691100
$ws.Range("H49").Value = 647000
$ws.Range("I49").Value = 617500
$ws.Range("J49").Value = 600900
$ws.Range("K49").Value = 607600
$ws.Range("D52").Value = 370000
$ws.Range("E52").Value = 445600
$ws.Range("F52").Value = 431800
$ws.Range("G52").Value = 410100
$ws.Range("H52").Value = 387800
$ws.Range("I52").Value = 362200
$ws.Range("J52").Value = 345500
$ws.Range("K52").Value = 345000
$ws.Range("D54").Value = 69843300
$ws.Range("E54").Value = 67396300
$ws.Range("F54").Value = 66198000
$ws.Range("G54").Value = 66211100
$ws.Range("H54").Value = 64714600
$ws.Range("I54").Value = 63449300
$ws.Range("J54").Value = 59398400
$ws.Range("K54").Value = 57315500
$ws.Range("D59").Value = 16400
$ws.Range("E59").Value = 17800
$ws.Range("F59").Value = 11900
$ws.Range("G59").Value = 14700
$ws.Range("H59").Value = 25400
$ws.Range("I59").Value = 28200
$ws.Range("J59").Value = 17500
$ws.Range("K59").Value = 9500
$ws.Range("D61").Value = 3679200
$ws.Range("E61").Value = 3535000
$ws.Range("F61").Value = 3834000
$ws.Range("G61").Value = 3569300
$ws.Range("H61").Value = 3550700
$ws.Range("I61").Value = 3083100
$ws.Range("J61").Value = 3100100
$ws.Range("K61").Value = 2523900
$ws.Range("D62").Value = 245200
$ws.Range("E62").Value = 293800
$ws.Range("F62").Value = 318700
$ws.Range("G62").Value = 291400
$ws.Range("H62").Value = 330600
$ws.Range("I62").Value = 291700
$ws.Range("J62").Value = 264500
$ws.Range("K62").Value = 283400
$ws.Range("D66").Value = 64035600
$ws.Range("E66").Value = 61607600
$ws.Range("F66").Value = 60671900
$ws.Range("G66").Value = 60883300
$ws.Range("H66").Value = 59421900
$ws.Range("I66").Value = 58325600
$ws.Range("J66").Value = 54691200
$ws.Range("K66").Value = 52568900
$ws.Range("D72").Value = 4147000
$ws.Range("E72").Value = 4243000
$ws.Range("F72").Value = 3961500
$ws.Range("G72").Value = 3815500
$ws.Range("H72").Value = 3698100
$ws.Range("I72").Value = 3527800
$ws.Range("J72").Value = 3311600
$ws.Range("K72").Value = 3525200
$ws.Range("D76").Value = 5807700
$ws.Range("E76").Value = 5788800
$ws.Range("F76").Value = 5526100
$ws.Range("G76").Value = 5327800
$ws.Range("H76").Value = 5292700
$ws.Range("I76").Value = 5123600
$ws.Range("J76").Value = 4707300
$ws.Range("K76").Value = 4746600
$ws.Range("D81").Value = 456300
$ws.Range("E81").Value = 437100
$ws.Range("F81").Value = 359400
$ws.Range("G81").Value = 321500
$ws.Range("H81").Value = 373000
$ws.Range("I81").Value = 369600
$ws.Range("J81").Value = 365200
$ws.Range("K81").Value = 368100
$ws.Range("D89").Value = 34800
$ws.Range("E89").Value = 462500
$ws.Range("F89").Value = -242200
$ws.Range("I89").Value = -288800
$ws.Range("J89").Value = 594000
$ws.Range("K89").Value = 740600
$ws.Range("D94").Value = -176600
$ws.Range("E94").Value = -176600
$ws.Range("F94").Value = -242800
$ws.Range("I94").Value = 152200
$ws.Range("J94").Value = 44300
$ws.Range("K94").Value = -494300
$ws.Range("D100").Value = -60500
$ws.Range("E100").Value = -481600
$ws.Range("F100").Value = 74100
$ws.Range("I100").Value = -203000
$ws.Range("J100").Value = 463600
$ws.Range("K100").Value = -203300
$ws.Range("D101").Value = -24300
$ws.Range("E101").Value = 7300
$ws.Range("F101").Value = -14900
$ws.Range("I101").Value = -15300
$ws.Range("J101").Value = -5300
$ws.Range("K101").Value = 1300
$ws.Range("D102").Value = -226500
$ws.Range("E102").Value = -188300
$ws.Range("F102").Value = -425900
$ws.Range("I102").Value = -354900
$ws.Range("J102").Value = 1096600
$ws.Range("K102").Value = 44300